$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "model_8_7_2"
$ws.Range("B2").Value = -0.3481427032956579
$ws.Range("C2").Value = -3.923842310695965
$ws.Range("D2").Value = 0.5053576727439146
$ws.Range("E2").Value = -1.555321257978936
$ws.Range("F2").Value = 1.491996765136719
$ws.Range("G2").Value = 6.318493366241455
$ws.Range("H2").Value = 0.8193596005439758
$ws.Range("I2").Value = 3.73066520690918
$ws.Range("A3").Value = "model_8_7_3"
$ws.Range("B3").Value = -0.2734161910119668
$ws.Range("C3").Value = -3.861244881459479
$ws.Range("D3").Value = 0.592954153361462
$ws.Range("E3").Value = -1.479419678856817
$ws.Range("F3").Value = 1.40929639339447
$ws.Range("G3").Value = 6.238165378570557
$ws.Range("H3").Value = 0.6742587089538574
$ws.Range("I3").Value = 3.619852066040039
$ws.Range("A4").Value = "model_8_7_1"
$ws.Range("B4").Value = -0.26143623276971
$ws.Range("C4").Value = -3.198286719001937
$ws.Range("D4").Value = 0.431631975987819
$ws.Range("E4").Value = -1.257059333967478
$ws.Range("F4").Value = 1.396038174629211
$ws.Range("G4").Value = 5.387428283691406
$ws.Range("H4").Value = 0.9414839744567871
$ws.Range("I4").Value = 3.295214653015137
$ws.Range("A5").Value = "model_8_7_0"
$ws.Range("B5").Value = -0.116232761665783
$ws.Range("C5").Value = -1.562819590997329
$ws.Range("D5").Value = 0.0630880305856486
$ws.Range("E5").Value = -0.6928051578191672
$ws.Range("F5").Value = 1.235340714454651
$ws.Range("G5").Value = 3.288723945617676
$ws.Range("H5").Value = 1.551965355873108
$ws.Range("I5").Value = 2.471426725387573
$ws.Range("A6").Value = "model_8_7_4"
$ws.Range("B6").Value = 0.04647971175569909
$ws.Range("C6").Value = -2.547662252658538
$ws.Range("D6").Value = 0.805721982135158
$ws.Range("E6").Value = -0.7545830994716505
$ws.Range("F6").Value = 1.05526602268219
$ws.Range("G6").Value = 4.552517890930176
$ws.Range("H6").Value = 0.3218154907226562
$ws.Range("I6").Value = 2.561619997024536
$ws.Range("A7").Value = "model_8_7_5"
$ws.Range("B7").Value = 0.33972322169777
$ws.Range("C7").Value = -1.128130688212573
$ws.Range("D7").Value = 0.8500733804426995
$ws.Range("E7").Value = -0.07033344762694993
$ws.Range("F7").Value = 0.7307318449020386
$ws.Range("G7").Value = 2.730912208557129
$ws.Range("H7").Value = 0.2483487725257874
$ws.Range("I7").Value = 1.56264328956604
$ws.Range("A8").Value = "model_8_7_6"
$ws.Range("B8").Value = 0.5156152873789634
$ws.Range("C8").Value = -0.4178553669616432
$ws.Range("D8").Value = 0.6085465405639059
$ws.Range("E8").Value = 0.1312188278450134
$ws.Range("F8").Value = 0.5360711216926575
$ws.Range("G8").Value = 1.819455027580261
$ws.Range("H8").Value = 0.6484304666519165
$ws.Range("I8").Value = 1.268385291099548
$ws.Range("A9").Value = "model_8_7_7"
$ws.Range("B9").Value = 0.5772594360779271
$ws.Range("C9").Value = -0.144954906624468
$ws.Range("D9").Value = 0.5345838137011412
$ws.Range("E9").Value = 0.2187178136540027
$ws.Range("F9").Value = 0.4678492248058319
$ws.Range("G9").Value = 1.469257116317749
$ws.Range("H9").Value = 0.7709474563598633
$ws.Range("I9").Value = 1.140640139579773
$ws.Range("A10").Value = "model_8_7_8"
$ws.Range("B10").Value = 0.5949153081662965
$ws.Range("C10").Value = 0.03359006000956044
$ws.Range("D10").Value = 0.4367195442608354
$ws.Range("E10").Value = 0.2495473059154244
$ws.Range("F10").Value = 0.4483094215393066
$ws.Range("G10").Value = 1.240140318870544
$ws.Range("H10").Value = 0.9330565333366394
$ws.Range("I10").Value = 1.095630407333374
$ws.Range("A11").Value = "model_8_7_19"
$ws.Range("B11").Value = 0.6355903165920409
$ws.Range("C11").Value = 0.4881443171259162
$ws.Range("D11").Value = -0.1056101476206188
$ws.Range("E11").Value = 0.1715002546677209
$ws.Range("F11").Value = 0.4032941460609436
$ws.Range("G11").Value = 0.6568359732627869
$ws.Range("H11").Value = 1.831408739089966
$ws.Range("I11").Value = 1.209576010704041
$ws.Range("A12").Value = "model_8_7_23"
$ws.Range("B12").Value = 0.6362833648054876
$ws.Range("C12").Value = 0.492739863790121
$ws.Range("D12").Value = -0.1063917963680419
$ws.Range("E12").Value = 0.1732220219965886
$ws.Range("F12").Value = 0.4025271236896515
$ws.Range("G12").Value = 0.6509387493133545
$ws.Range("H12").Value = 1.832703471183777
$ws.Range("I12").Value = 1.207062244415283
$ws.Range("A13").Value = "model_8_7_22"
$ws.Range("B13").Value = 0.6364157646142306
$ws.Range("C13").Value = 0.4929235002114769
$ws.Range("D13").Value = -0.105924202991202
$ws.Range("E13").Value = 0.1735527823833495
$ws.Range("F13").Value = 0.4023806154727936
$ws.Range("G13").Value = 0.6507031321525574
$ws.Range("H13").Value = 1.831928849220276
$ws.Range("I13").Value = 1.206579327583313
$ws.Range("A14").Value = "model_8_7_24"
$ws.Range("B14").Value = 0.6364231588087779
$ws.Range("C14").Value = 0.4925244564204331
$ws.Range("D14").Value = -0.1055470462020018
$ws.Range("E14").Value = 0.1735700971571591
$ws.Range("F14").Value = 0.402372419834137
$ws.Range("G14").Value = 0.6512151956558228
$ws.Range("H14").Value = 1.83130419254303
$ws.Range("I14").Value = 1.206553936004639
$ws.Range("A15").Value = "model_8_7_21"
$ws.Range("B15").Value = 0.6365134228169216
$ws.Range("C15").Value = 0.4934733890573527
$ws.Range("D15").Value = -0.1059527728253866
$ws.Range("E15").Value = 0.1737972231503985
$ws.Range("F15").Value = 0.4022725224494934
$ws.Range("G15").Value = 0.6499974727630615
$ws.Range("H15").Value = 1.831976294517517
$ws.Range("I15").Value = 1.206222295761108
$ws.Range("A16").Value = "model_8_7_20"
$ws.Range("B16").Value = 0.6365574980818984
$ws.Range("C16").Value = 0.4936605675517753
$ws.Range("D16").Value = -0.105909474542015
$ws.Range("E16").Value = 0.1739081363194395
$ws.Range("F16").Value = 0.4022237360477448
$ws.Range("G16").Value = 0.6497572660446167
$ws.Range("H16").Value = 1.831904649734497
$ws.Range("I16").Value = 1.206060528755188
$ws.Range("A17").Value = "model_8_7_18"
$ws.Range("B17").Value = 0.6393256978476642
$ws.Range("C17").Value = 0.5049491437359589
$ws.Range("D17").Value = -0.102769096955885
$ws.Range("E17").Value = 0.1808359662079231
$ws.Range("F17").Value = 0.3991601765155792
$ws.Range("G17").Value = 0.6352713108062744
$ws.Range("H17").Value = 1.82670259475708
$ws.Range("I17").Value = 1.195946216583252
$ws.Range("A18").Value = "model_8_7_17"
$ws.Range("B18").Value = 0.6393442131807334
$ws.Range("C18").Value = 0.5054345157930575
$ws.Range("D18").Value = -0.103096012859671
$ws.Range("E18").Value = 0.1808886044980023
$ws.Range("F18").Value = 0.3991396725177765
$ws.Range("G18").Value = 0.6346484422683716
$ws.Range("H18").Value = 1.827244162559509
$ws.Range("I18").Value = 1.195869207382202
$ws.Range("A19").Value = "model_8_7_16"
$ws.Range("B19").Value = 0.6398139162598435
$ws.Range("C19").Value = 0.5077590900021137
$ws.Range("D19").Value = -0.1029110465529197
$ws.Range("E19").Value = 0.182069849996419
$ws.Range("F19").Value = 0.3986198604106903
$ws.Range("G19").Value = 0.6316654682159424
$ws.Range("H19").Value = 1.826937794685364
$ws.Range("I19").Value = 1.194144606590271
$ws.Range("A20").Value = "model_8_7_15"
$ws.Range("B20").Value = 0.6424927276945516
$ws.Range("C20").Value = 0.5112652818762047
$ws.Range("D20").Value = -0.09329708598506414
$ws.Range("E20").Value = 0.1888326176626701
$ws.Range("F20").Value = 0.3956551849842072
$ws.Range("G20").Value = 0.6271660923957825
$ws.Range("H20").Value = 1.811012506484985
$ws.Range("I20").Value = 1.184271335601807
$ws.Range("A21").Value = "model_8_7_14"
$ws.Range("B21").Value = 0.6454733502969633
$ws.Range("C21").Value = 0.5112364296772032
$ws.Range("D21").Value = -0.0785588794435228
$ws.Range("E21").Value = 0.1966894596775595
$ws.Range("F21").Value = 0.3923565447330475
$ws.Range("G21").Value = 0.6272031664848328
$ws.Range("H21").Value = 1.786599159240723
$ws.Range("I21").Value = 1.172800660133362
$ws.Range("A22").Value = "model_8_7_12"
$ws.Range("B22").Value = 0.6470163300712926
$ws.Range("C22").Value = 0.4434181808624107
$ws.Range("D22").Value = -0.0007133490040469592
$ws.Range("E22").Value = 0.2066951781528357
$ws.Range("F22").Value = 0.3906489014625549
$ws.Range("G22").Value = 0.7142305374145508
$ws.Range("H22").Value = 1.657650470733643
$ws.Range("I22").Value = 1.15819263458252
$ws.Range("A23").Value = "model_8_7_13"
$ws.Range("B23").Value = 0.65193794622097
$ws.Range("C23").Value = 0.4924211624904572
$ws.Range("D23").Value = -0.02965048215907262
$ws.Range("E23").Value = 0.2140461683281487
$ws.Range("F23").Value = 0.3852021396160126
$ws.Range("G23").Value = 0.6513477563858032
$ws.Range("H23").Value = 1.705584049224854
$ws.Range("I23").Value = 1.147460579872131
$ws.Range("A24").Value = "model_8_7_11"
$ws.Range("B24").Value = 0.6724227544034771
$ws.Range("C24").Value = 0.4318192861271871
$ws.Range("D24").Value = 0.1514664257265504
$ws.Range("E24").Value = 0.2825498998374109
$ws.Range("F24").Value = 0.362531453371048
$ws.Range("G24").Value = 0.7291147708892822
$ws.Range("H24").Value = 1.405569314956665
$ws.Range("I24").Value = 1.047447919845581
$ws.Range("A25").Value = "model_8_7_9"
$ws.Range("B25").Value = 0.6897162846891278
$ws.Range("C25").Value = 0.3520840029069305
$ws.Range("D25").Value = 0.5578149274031472
$ws.Range("E25").Value = 0.4624078936717247
$ws.Range("F25").Value = 0.3433926403522491
$ws.Range("G25").Value = 0.8314345479011536
$ws.Range("H25").Value = 0.7324658036231995
$ws.Range("I25").Value = 0.7848625183105469
$ws.Range("A26").Value = "model_8_7_10"
$ws.Range("B26").Value = 0.7035704685378332
$ws.Range("C26").Value = 0.5484584081870258
$ws.Range("D26").Value = 0.2962096172653599
$ws.Range("E26").Value = 0.4141103309631925
$ws.Range("F26").Value = 0.3280601501464844
$ws.Range("G26").Value = 0.579438328742981
$ws.Range("H26").Value = 1.165806889533997
$ws.Range("I26").Value = 0.8553750514984131
